$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) ---
$ws.Range("A1").Value = "Reference"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "JLCPCB Part #"
$ws.Range("D1").Value = "Price"
$ws.Range("E1").Value = "Qty"

# --- Row 2 (C1, C2 / 100nF) ---
$ws.Range("A2").Value = ">  C1, C2"
$ws.Range("B2").Value = "100nF"
$ws.Range("C2").Value = "C254109"
$ws.Range("D2").Value = "CC2A104ZC1ID3F7C30MF"
$ws.Range("E2").Value = "`$0.0229"

# --- Row 3 (EEPROM1) ---
$ws.Range("A3").Value = "    EEPROM1"
$ws.Range("B3").Value = "AT24C256C-SSHL-T EEPROM"
$ws.Range("C3").Value = "C6482"
$ws.Range("D3").Value = "AT24C256C-SSHL-T"
$ws.Range("E3").Value = "`$1.5556"

# --- Row 4 (MCU1) ---
$ws.Range("A4").Value = "    MCU1"
$ws.Range("B4").Value = "STM32F030C8Tx"
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = "`$0.00"

# --- Row 5 (OP-AMP1) ---
$ws.Range("A5").Value = "    OP-AMP1"
$ws.Range("B5").Value = "LM358DT"
$ws.Range("C5").Value = "C9418"
$ws.Range("D5").Value = "LM358DT"
$ws.Range("E5").Value = "`$0.1043"

# --- Row 6 (R1, R2 / 33k) ---
$ws.Range("A6").Value = ">  R1, R2"
$ws.Range("B6").Value = "33k"
$ws.Range("C6").Value = "C425317"
$ws.Range("D6").Value = "4D02WGF3302TCE"
$ws.Range("E6").Value = "`$0.006"

# --- Resize the table + enable totals row, then set total label text ---
$lo = $ws.ListObjects.Item(1)
$lo.ShowTotals = $true

$totalsRange = $lo.TotalsRowRange
$totalsRange.Cells.Item(1,5).Value = "Total = `$1.6892"

$col = $lo.ListColumns.Item(5)
$col.TotalsCalculation = 9

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 23.81640625
$ws.Columns.Item(5).ColumnWidth = 21.90625

# --- Selection ---
$ws.Range("D10").Select()
